$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.400.28"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "3.083.10"
$ws.Range("E3").Value = "  +4.16%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'580.86"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("D6").Value = "'167.11"
$ws.Range("E6").Value = "  +4.19%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.079.32"
$ws.Range("E8").Value = "  +4.09%  "

$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  +1.43%  "

$ws.Range("D10").Value = "'6.67"
$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("E11").Value = "  +1.27%  "

$ws.Range("E12").Value = "  +5.82%  "

$ws.Range("D13").Value = "'0.0000248"

$ws.Range("D14").Value = "'36.60"
$ws.Range("E14").Value = "  +7.09%  "

$ws.Range("D16").Value = "3.595.20"
$ws.Range("E16").Value = "  +4.10%  "

$ws.Range("D17").Value = "66.378.26"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").Value = "'7.16"
$ws.Range("E18").Value = "  +4.08%  "

$ws.Range("D19").Value = "3.083.82"
$ws.Range("E19").Value = "  +4.20%  "

$ws.Range("D20").Value = "'15.98"
$ws.Range("E20").Value = "  +16.79%  "

$ws.Range("D21").Value = "'462.42"
$ws.Range("E21").Value = "  +3.79%  "

$ws.Range("E22").Value = "  +5.90%  "

$ws.Range("E23").Value = "  +4.46%  "

$ws.Range("D24").Value = "'83.09"
$ws.Range("E24").Value = "  +1.50%  "

$ws.Range("D25").Value = "'12.76"
$ws.Range("E25").Value = "  +5.04%  "

$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  +0.95%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").Value = "'8.02"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").Value = "'2.40"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("E31").Value = "  +3.29%  "

$ws.Range("E32").Value = "  +3.94%  "

$ws.Range("D33").Value = "'28.39"
$ws.Range("E33").Value = "  +4.89%  "

$ws.Range("E34").Value = "  +5.65%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "'0.997"

$ws.Range("D37").Value = "'5.87"
$ws.Range("E37").Value = "  +3.59%  "

$ws.Range("D38").Value = "'48.21"
$ws.Range("E38").Value = "  +12.02%  "

$ws.Range("D39").Value = "'50.09"
$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.04"
$ws.Range("E40").Value = "  +2.86%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.312"
$ws.Range("E41").Value = "  +4.46%  "

$ws.Range("E42").Value = "  +2.58%  "

$ws.Range("D43").Value = "'2.88"
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("D44").Value = "'8.63"
$ws.Range("E44").Value = "  +3.63%  "

$ws.Range("E45").Value = "  +2.11%  "

$ws.Range("D46").Value = "'384.23"
$ws.Range("E46").Value = "  +1.15%  "

$ws.Range("D47").Value = "2.773.30"
$ws.Range("E47").Value = "  +2.45%  "

$ws.Range("D48").Value = "'134.82"
$ws.Range("E48").Value = "  +2.75%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").Value = "'24.48"
$ws.Range("E50").Value = "  +6.49%  "

$ws.Range("E51").Value = "  +4.58%  "
